$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.761.80'
$ws.Range("E2").Value = '  -0.25%  '
$ws.Range("D3").Value = '2.307.45'
$ws.Range("E3").Value = '  +3.46%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '270.60'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +0.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.13'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +5.68%  '
$ws.Range("E7").Value = '  +0.94%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.623'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +2.44%  '
$ws.Range("E10").Value = '  -3.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0937'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +1.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.07'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +6.44%  '
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("D14").Value = '2.654.02'
$ws.Range("E14").Value = '  +3.47%  '
$ws.Range("E15").Value = '  +3.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.849'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +6.87%  '
$ws.Range("D17").Value = '2.282.00'
$ws.Range("E17").Value = '  +1.33%  '
$ws.Range("D18").Value = '43.703.77'
$ws.Range("E18").Value = '  -0.37%  '
$ws.Range("E19").Value = '  +2.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.26'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +3.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.45'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +1.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '240.27'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +3.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.29'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -5.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.73'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +9.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.36'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +4.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.50'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -4.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.37'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +5.49%  '
$ws.Range("E29").Value = '  -5.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.99'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -3.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.41'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  +8.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '171.67'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -2.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0897'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -1.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.58'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +2.79%  '
$ws.Range("E35").Value = '  +1.61%  '
$ws.Range("E36").Value = '  +0.13%  '
$ws.Range("E37").Value = '  +2.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0353'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -1.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.43'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +2.34%  '
$ws.Range("E40").Value = '  +14.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.29'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +8.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.17'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -2.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.31'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +15.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.44'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +1.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.47'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -6.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.89'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +6.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.103'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +3.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '100.28'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -1.01%  '
$ws.Range("E49").Value = '  -2.04%  '
$ws.Range("D50").Value = '2.531.54'
$ws.Range("E50").Value = '  +3.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.425'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -3.42%  '
